$d = $word.ActiveDocument

function TestPos($pos) {
    $target = $d.Range($pos, $pos)
    $name = "_TestBm"
    if ($d.Bookmarks.Exists($name)) {
        $d.Bookmarks.Item($name).Delete()
    }
    $bm = $d.Bookmarks.Add($name, $target)
    $bm2 = $d.Bookmarks.Item($name)
    Write-Output ("pos=" + $pos + " -> bm start=" + $bm2.Start + " end=" + $bm2.End)
}

for ($p = 1788; $p -le 1800; $p++) {
    TestPos $p
}
